$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '50.873.89'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -16.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.230.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -23.30%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '424.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -22.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.442'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -19.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.228.36'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -23.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -16.93%  '
$ws.Range('E11').Value = '  -24.74%  '
$ws.Range('E12').Value = '  -7.12%  '
$ws.Range('E13').Value = '  -19.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.617.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -23.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '50.962.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -16.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '17.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -21.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.244.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -22.52%  '
$ws.Range('E18').Value = '  -22.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -24.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '282.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -20.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.991'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -29.53%  '
$ws.Range('E24').Value = '  -27.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '51.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -20.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.329.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -23.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.348'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -22.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.131'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -26.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -18.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '141.76'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0₃0591'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -31.02%  '
$ws.Range('E34').Value = '  -18.78%  '
$ws.Range('E35').Value = '  -24.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.43'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -20.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.995'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.748'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -24.83%  '
$ws.Range('E39').Value = '  -27.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '31.06'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -17.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.926'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -22.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.12'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.537'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -17.56%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.93'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -21.07%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0470'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -19.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.811.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -20.96%  '
$ws.Range('E47').Value = '  -27.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0193'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -18.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0770'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -16.00%  '
$ws.Range('E50').Value = '  -5.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '14.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -27.40%  '
